$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("B1").Value = "sparsity_necessary"
$ws.Range("C1").Value = "necessary explanation rate"

# Update "1-best" -> "1-delta" for CoDy rows (2-16)
$ws.Range("E2:E16").Value = "1-delta"

# Update "recent" -> "temporal" for Greedy rows (41-48)
$ws.Range("E41:E48").Value = "temporal"
